$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.681.35'
$ws.Range("E2").Value = '  +6.16%  '

$ws.Range("D3").Value = '2.055.08'
$ws.Range("E3").Value = '  +3.11%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.99'
$ws.Range("E5").Value = '  +5.10%  '

$ws.Range("E6").Value = '  +2.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.26'
$ws.Range("E7").Value = '  +17.85%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.377'
$ws.Range("E9").Value = '  +6.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.44'
$ws.Range("E10").Value = '  +0.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0757'
$ws.Range("E11").Value = '  +4.43%  '

$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.35'
$ws.Range("E13").Value = '  +7.43%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.913'
$ws.Range("E14").Value = '  +2.11%  '

$ws.Range("D15").Value = '2.355.18'
$ws.Range("E15").Value = '  +3.30%  '

$ws.Range("E16").Value = '  +8.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.88'
$ws.Range("E17").Value = '  +22.65%  '

$ws.Range("D18").Value = '2.074.59'
$ws.Range("E18").Value = '  +4.18%  '

$ws.Range("D19").Value = '37.565.57'
$ws.Range("E19").Value = '  +6.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.23'
$ws.Range("E20").Value = '  +4.93%  '

$ws.Range("D21").Value = '0.0₃0879'
$ws.Range("E21").Value = '  +5.47%  '

$ws.Range("E22").Value = '  +6.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.73'
$ws.Range("E23").Value = '  +2.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.79'
$ws.Range("E24").Value = '  +23.74%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("E26").Value = '  +4.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.62'
$ws.Range("E27").Value = '  +5.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.54'
$ws.Range("E28").Value = '  +1.56%  '

$ws.Range("E29").Value = '  +2.54%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.24'
$ws.Range("E30").Value = '  +9.88%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("E31").Value = '  +8.20%  '

$ws.Range("E32").Value = '  +2.85%  '

$ws.Range("E33").Value = '  +25.26%  '

$ws.Range("E34").Value = '  +12.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0617'
$ws.Range("E35").Value = '  +5.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("E36").Value = '  +10.27%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.02'
$ws.Range("E38").Value = '  +23.74%  '

$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("E40").Value = '  +18.20%  '

$ws.Range("E41").Value = '  +5.07%  '

$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.93'
$ws.Range("E42").Value = '  +4.14%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +5.43%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.17'
$ws.Range("E44").Value = '  +10.47%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.14'
$ws.Range("E45").Value = '  +5.97%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.13'
$ws.Range("E46").Value = '  +11.27%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").Value = '  +19.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '95.71'
$ws.Range("E48").Value = '  +6.18%  '

$ws.Range("D49").Value = '1.431.40'
$ws.Range("E49").Value = '  +4.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  +1.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.65'
$ws.Range("E51").Value = '  +5.47%  '
